$p = $ppt.ActivePresentation

# --- 1. Reposition shapes on Slide 2 ("Cross Join" explanation slide) ---
$s2 = $p.Slides.Item(2)

# "Rectangle 35" textbox: move up (only Y changes)
$rect = $s2.Shapes.Item("Rectangle 35")
$rect.Left = 380166 / 12700
$rect.Top = 1142984 / 12700

# Second picture (diagram image) on the slide: reposition X and Y
$pic = $s2.Shapes.Item(4)
$pic.Left = 1345435 / 12700
$pic.Top = 2205058 / 12700

# --- 2. Update the cached "datetimeFigureOut" field text (6/18/2020 -> 6/22/2020) ---
# on every slide layout's Date Placeholder ...
$master = $p.Designs.Item(1).SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {}
        if ($isDatePh -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "6/18/2020") {
                $shp.TextFrame.TextRange.Text = "6/22/2020"
            }
        }
    }
}

# ... and on the Notes Master's Date Placeholder
$notesMaster = $p.NotesMaster
for ($si = 1; $si -le $notesMaster.Shapes.Count; $si++) {
    $shp = $notesMaster.Shapes.Item($si)
    $isDatePh = $false
    try {
        if ($shp.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
    } catch {}
    if ($isDatePh -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "6/18/2020") {
            $shp.TextFrame.TextRange.Text = "6/22/2020"
        }
    }
}
